$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D price values remain stored as text (they are mixed
# thousand-separated / plain decimal strings in the source data), not
# auto-converted to numbers, by forcing a Text format before writing and
# then resetting the style back to Normal so no stray style id is left
# behind on the cells.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '29.970.35'
$ws.Range("E2").Value = '  +0.64%  '

$ws.Range("D3").Value = '1.909.55'
$ws.Range("E3").Value = '  +1.05%  '

$ws.Range("D4").Value = '0.9991'
$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").Value = '0.8175'
$ws.Range("E5").Value = '  +8.33%  '

$ws.Range("D6").Value = '241.57'
$ws.Range("E6").Value = '  +1.00%  '

$ws.Range("D7").Value = '0.9973'
$ws.Range("E7").Value = '  -0.24%  '

$ws.Range("D8").Value = '0.3163'
$ws.Range("E8").Value = '  +4.14%  '

$ws.Range("D9").Value = '26.64'
$ws.Range("E9").Value = '  +4.70%  '

$ws.Range("D10").Value = '0.07008'
$ws.Range("E10").Value = '  +2.96%  '

$ws.Range("D11").Value = '0.08011'
$ws.Range("E11").Value = '  +0.97%  '

$ws.Range("D12").Value = '0.7481'
$ws.Range("E12").Value = '  +0.83%  '

$ws.Range("D13").Value = '1.907.90'
$ws.Range("E13").Value = '  +1.03%  '

$ws.Range("D14").Value = '5.195'
$ws.Range("E14").Value = '  +1.12%  '

$ws.Range("D15").Value = '92.57'
$ws.Range("E15").Value = '  +2.34%  '

$ws.Range("D16").Value = '29.973.67'
$ws.Range("E16").Value = '  +0.67%  '

$ws.Range("D17").Value = '14.10'
$ws.Range("E17").Value = '  +1.63%  '

$ws.Range("D18").Value = '5.900'
$ws.Range("E18").Value = '  -0.65%  '

$ws.Range("D19").Value = '245.70'
$ws.Range("E19").Value = '  +1.79%  '

$ws.Range("D20").Value = '0.000007778'
$ws.Range("E20").Value = '  +1.54%  '

$ws.Range("D21").Value = '1.004'
$ws.Range("E21").Value = '  +0.45%  '

$ws.Range("D22").Value = '2.153.25'
$ws.Range("E22").Value = '  +0.53%  '

$ws.Range("D23").Value = '0.9990'
$ws.Range("E23").Value = '  -0.10%  '

$ws.Range("D24").Value = '6.970'
$ws.Range("E24").Value = '  +0.88%  '

$ws.Range("D25").Value = '0.1586'
$ws.Range("E25").Value = '  +24.49%  '

$ws.Range("D26").Value = '167.88'
$ws.Range("E26").Value = '  +1.40%  '

$ws.Range("D27").Value = '9.233'
$ws.Range("E27").Value = '  +0.51%  '

$ws.Range("D28").Value = '18.88'
$ws.Range("E28").Value = '  +1.54%  '

$ws.Range("D29").Value = '2.086'
$ws.Range("E29").Value = '  +3.53%  '

$ws.Range("D30").Value = '1.364'
$ws.Range("E30").Value = '  -1.31%  '

$ws.Range("E31").Value = '  -0.10%  '

$ws.Range("D32").Value = '4.313'
$ws.Range("E32").Value = '  +1.90%  '

$ws.Range("D33").Value = '4.099'
$ws.Range("E33").Value = '  +2.21%  '

$ws.Range("D34").Value = '0.05533'
$ws.Range("E34").Value = '  +6.70%  '

$ws.Range("D35").Value = '1.272'
$ws.Range("E35").Value = '  +1.93%  '

$ws.Range("D36").Value = '0.7353'
$ws.Range("E36").Value = '  +1.41%  '

$ws.Range("D37").Value = '2.694'
$ws.Range("E37").Value = '  -0.38%  '

$ws.Range("D38").Value = '0.01923'
$ws.Range("E38").Value = '  +1.14%  '

$ws.Range("D39").Value = '2.787'
$ws.Range("E39").Value = '  +0.74%  '

$ws.Range("D40").Value = '0.4424'
$ws.Range("E40").Value = '  +1.05%  '

$ws.Range("D41").Value = '72.45'
$ws.Range("E41").Value = '  +1.81%  '

$ws.Range("D42").Value = '5.986'
$ws.Range("E42").Value = '  -2.21%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '1.900'
$ws.Range("E43").Value = '  +1.03%  '

$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").Value = '0.9962'
$ws.Range("E44").Value = '  -0.34%  '

$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").Value = '0.8361'
$ws.Range("E45").Value = '  +1.41%  '

$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").Value = '7.591'
$ws.Range("E46").Value = '  +0.17%  '

$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").Value = '100.82'
$ws.Range("E47").Value = '  +1.20%  '

$ws.Range("D48").Value = '9.734'
$ws.Range("E48").Value = '  +0.30%  '

$ws.Range("D49").Value = '986.08'
$ws.Range("E49").Value = '  +9.81%  '

$ws.Range("D50").Value = '2.059.28'
$ws.Range("E50").Value = '  +0.92%  '

$ws.Range("E51").Value = '  +0.91%  '

# Restore the default "Normal" style on the price column now that the
# text values are safely written, matching the original formatting.
$priceRange.Style = "Normal"
